$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E stay formatted as text so numeric-looking
# strings (e.g. "98.608.28", "1.00", "  -0.65%  ") are preserved verbatim
# instead of being re-interpreted as numbers/dates by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.608.28"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "3.344.94"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "258.32"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").Value = "664.60"
$ws.Range("E6").Value = "  +5.69%  "

$ws.Range("E7").Value = "  +12.40%  "

$ws.Range("D8").Value = "0.461"
$ws.Range("E8").Value = "  +17.47%  "

$ws.Range("E9").Value = "  +27.25%  "

$ws.Range("D11").Value = "3.343.26"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("E12").Value = "  +5.27%  "

$ws.Range("D13").Value = "42.90"
$ws.Range("E13").Value = "  +19.54%  "

$ws.Range("D14").Value = "0.0000268"
$ws.Range("E14").Value = "  +8.36%  "

$ws.Range("D15").Value = "99.227.78"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "3.965.65"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "5.65"
$ws.Range("E17").Value = "  +2.78%  "

$ws.Range("D18").Value = "3.338.23"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  +24.49%  "

$ws.Range("D20").Value = "17.00"
$ws.Range("E20").Value = "  +11.40%  "

$ws.Range("D21").Value = "3.62"
$ws.Range("E21").Value = "  +1.98%  "

$ws.Range("D22").Value = "529.26"
$ws.Range("E22").Value = "  +7.29%  "

$ws.Range("D23").Value = "10.50"
$ws.Range("E23").Value = "  +12.19%  "

$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "0.0000211"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "0.437"
$ws.Range("E25").Value = "  +56.62%  "

$ws.Range("D26").Value = "6.38"
$ws.Range("E26").Value = "  +12.89%  "

$ws.Range("D27").Value = "101.40"
$ws.Range("E27").Value = "  +14.88%  "

$ws.Range("D28").Value = "12.83"
$ws.Range("E28").Value = "  +7.68%  "

$ws.Range("D29").Value = "3.519.12"
$ws.Range("E29").Value = "  -0.80%  "

$ws.Range("D30").Value = "0.149"
$ws.Range("E30").Value = "  +15.82%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "11.15"
$ws.Range("E32").Value = "  +16.93%  "

$ws.Range("D33").Value = "0.191"
$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "29.33"
$ws.Range("E35").Value = "  +4.97%  "

$ws.Range("D36").Value = "0.539"
$ws.Range("E36").Value = "  +17.45%  "

$ws.Range("D37").Value = "7.93"
$ws.Range("E37").Value = "  +8.26%  "

$ws.Range("E38").Value = "  +7.99%  "

$ws.Range("D39").Value = "0.160"
$ws.Range("E39").Value = "  +6.12%  "

$ws.Range("D40").Value = "533.61"

$ws.Range("D41").Value = "1.35"
$ws.Range("E41").Value = "  +6.98%  "

$ws.Range("D42").Value = "0.0455"
$ws.Range("E42").Value = "  +39.78%  "

$ws.Range("D43").Value = "24.71"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("E44").Value = "  -3.80%  "

$ws.Range("D45").Value = "0.827"
$ws.Range("E45").Value = "  +6.16%  "

$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +2.71%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +7.79%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "5.18"
$ws.Range("E49").Value = "  +11.96%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "7.95"
$ws.Range("E50").Value = "  +21.94%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "1.53"
$ws.Range("E51").Value = "  +12.15%  "
